$wb = $excel.ActiveWorkbook
$cmds = $wb.list_commands()
Write-Output $cmds
